$wb = $excel.ActiveWorkbook

# Remove the bold/border/center-top header style from row 1 (A1:N1) on every sheet
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Range("A1:N1").ClearFormats()
}

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H4").Value = 154.625
$ws.Range("I4").Value = 159.28572
$ws.Range("J4").Value = 122
$ws.Range("K4").Value = 159.28572
$ws.Range("L4").Value = 122
$ws.Range("M4").Value = -45.28572
$ws.Range("N4").Value = -350
$ws.Range("H41").Value = 312.72726
$ws.Range("I41").Value = 293.33334
$ws.Range("K41").Value = 293.33334
$ws.Range("M41").Value = 146.66666
$ws.Range("H51").Value = 2067.25
$ws.Range("I51").Value = 2089
$ws.Range("J51").Value = 2060
$ws.Range("K51").Value = 2089
$ws.Range("L51").Value = 2060
$ws.Range("M51").Value = -1605
$ws.Range("N51").Value = -3028
$ws.Range("H55").Value = 150
$ws.Range("I55").Value = 113.77778
$ws.Range("J55").Value = 215.2
$ws.Range("K55").Value = 113.77778
$ws.Range("L55").Value = 215.2
$ws.Range("M55").Value = 100.22222
$ws.Range("N55").Value = -643.2
$ws.Range("H64").Value = 12750
$ws.Range("J64").Value = 12750
$ws.Range("L64").Value = 12750
$ws.Range("N64").Value = -13246
$ws.Range("H67").Value = 12750
$ws.Range("J67").Value = 12750
$ws.Range("L67").Value = 12750
$ws.Range("N67").Value = -14466
$ws.Range("H82").Value = 1450
$ws.Range("I82").Value = 1450
$ws.Range("K82").Value = 4350
$ws.Range("M82").Value = -3944
$ws.Range("H85").Value = 1450
$ws.Range("I85").Value = 1450
$ws.Range("K85").Value = 4350
$ws.Range("M85").Value = -2946
$ws.Range("H131").Value = 43462.42
$ws.Range("I131").Value = 47327.566
$ws.Range("K131").Value = 141982.698
$ws.Range("M131").Value = -136942.698
$ws.Range("H135").Value = 1717.3125
$ws.Range("I135").Value = 1587.5834
$ws.Range("J135").Value = 2106.5
$ws.Range("K135").Value = 14288.2506
$ws.Range("L135").Value = 18958.5
$ws.Range("M135").Value = -11753.2506
$ws.Range("N135").Value = -24028.5
$ws.Range("H138").Value = 4239.1333
$ws.Range("I138").Value = 6284.75
$ws.Range("J138").Value = 3495.2727
$ws.Range("K138").Value = 18854.25
$ws.Range("L138").Value = 10485.8181
$ws.Range("M138").Value = -13714.25
$ws.Range("N138").Value = -20765.8181
$ws.Range("H141").Value = 2920.6956
$ws.Range("J141").Value = 5000
$ws.Range("L141").Value = 15000
$ws.Range("N141").Value = -25360

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 6670.087
$ws.Range("I32").Value = 5880.372
$ws.Range("K32").Value = 5880.372
$ws.Range("M32").Value = -5593.372
$ws.Range("H61").Value = 8663.75
$ws.Range("I61").Value = 7439.875
$ws.Range("K61").Value = 7439.875
$ws.Range("M61").Value = -7227.875
$ws.Range("H74").Value = 3703.68
$ws.Range("I74").Value = 2587.1875
$ws.Range("K74").Value = 2587.1875
$ws.Range("M74").Value = -1713.1875
$ws.Range("H77").Value = 3703.68
$ws.Range("I77").Value = 2587.1875
$ws.Range("K77").Value = 12935.9375
$ws.Range("M77").Value = -8567.9375
$ws.Range("H122").Value = 3121.7896
$ws.Range("I122").Value = 2565.7856
$ws.Range("J122").Value = 4678.6
$ws.Range("K122").Value = 7697.3568
$ws.Range("L122").Value = 14035.8
$ws.Range("M122").Value = -5247.3568
$ws.Range("N122").Value = -18935.8
$ws.Range("H132").Value = 5444.2856
$ws.Range("I132").Value = 4704.6816
$ws.Range("J132").Value = 8156.1665
$ws.Range("K132").Value = 14114.0448
$ws.Range("L132").Value = 24468.4995
$ws.Range("M132").Value = -11584.0448
$ws.Range("N132").Value = -29528.4995
$ws.Range("H136").Value = 8663.75
$ws.Range("I136").Value = 7439.875
$ws.Range("K136").Value = 22319.625
$ws.Range("M136").Value = -19769.625

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H20").Value = 3117.8333
$ws.Range("I20").Value = 3212.7778
$ws.Range("K20").Value = 3212.7778
$ws.Range("M20").Value = -2965.7778
$ws.Range("H86").Value = 13337980
$ws.Range("I86").Value = 5522.8887
$ws.Range("K86").Value = 5522.8887
$ws.Range("M86").Value = -4399.8887
$ws.Range("H89").Value = 13337980
$ws.Range("I89").Value = 5522.8887
$ws.Range("K89").Value = 27614.4435
$ws.Range("M89").Value = -21998.4435
$ws.Range("H94").Value = 1707.6562
$ws.Range("I94").Value = 1299.6072
$ws.Range("K94").Value = 1299.6072
$ws.Range("M94").Value = -848.6071999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H7").Value = 210.73077
$ws.Range("I7").Value = 140.5625
$ws.Range("J7").Value = 323
$ws.Range("K7").Value = 140.5625
$ws.Range("L7").Value = 323
$ws.Range("M7").Value = -27.5625
$ws.Range("N7").Value = -549
$ws.Range("H31").Value = 4173.875
$ws.Range("I31").Value = 3684.111
$ws.Range("J31").Value = 4574.591
$ws.Range("K31").Value = 3684.111
$ws.Range("L31").Value = 4574.591
$ws.Range("M31").Value = -3389.111
$ws.Range("N31").Value = -5164.591
$ws.Range("H34").Value = 4173.875
$ws.Range("I34").Value = 3684.111
$ws.Range("J34").Value = 4574.591
$ws.Range("K34").Value = 3684.111
$ws.Range("L34").Value = 4574.591
$ws.Range("M34").Value = -3482.111
$ws.Range("N34").Value = -4978.591
$ws.Range("H35").Value = 681
$ws.Range("I35").Value = 681
$ws.Range("K35").Value = 681
$ws.Range("M35").Value = -387
$ws.Range("H133").Value = 85100
$ws.Range("J133").Value = 85100
$ws.Range("L133").Value = 85100
$ws.Range("N133").Value = -90160
$ws.Range("H134").Value = 6719.2856
$ws.Range("I134").Value = 2759.6875
$ws.Range("K134").Value = 8279.0625
$ws.Range("M134").Value = -5744.0625

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H5").Value = 1740.7667
$ws.Range("I5").Value = 1610.5714
$ws.Range("J5").Value = 1854.6875
$ws.Range("K5").Value = 4831.7142
$ws.Range("L5").Value = 5564.0625
$ws.Range("M5").Value = -4719.7142
$ws.Range("N5").Value = -5788.0625
$ws.Range("H51").Value = 22495
$ws.Range("J51").Value = 30252.75
$ws.Range("L51").Value = 90758.25
$ws.Range("N51").Value = -91678.25
$ws.Range("J131").Value = 3037.2222
$ws.Range("L131").Value = 9111.6666
$ws.Range("N131").Value = -19191.6666
$ws.Range("H135").Value = 1740.7667
$ws.Range("I135").Value = 1610.5714
$ws.Range("J135").Value = 1854.6875
$ws.Range("K135").Value = 14495.1426
$ws.Range("L135").Value = 16692.1875
$ws.Range("M135").Value = -11960.1426
$ws.Range("N135").Value = -21762.1875
$ws.Range("H141").Value = 2749.5
$ws.Range("I141").Value = 2749.5
$ws.Range("K141").Value = 8248.5
$ws.Range("M141").Value = -3068.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H11").Value = 10500000
$ws.Range("I11").Value = 20000000
$ws.Range("J11").Value = 1000000
$ws.Range("K11").Value = 20000000
$ws.Range("L11").Value = 1000000
$ws.Range("M11").Value = -19999861
$ws.Range("N11").Value = -1000278
$ws.Range("H107").Value = 578.6667
$ws.Range("I107").Value = 118
$ws.Range("K107").Value = 118
$ws.Range("M107").Value = 1802
$ws.Range("H113").Value = 449344.34
$ws.Range("J113").Value = 5025
$ws.Range("L113").Value = 5025
$ws.Range("N113").Value = -9365
$ws.Range("H122").Value = 4901.909
$ws.Range("I122").Value = 5221.6
$ws.Range("J122").Value = 1705
$ws.Range("K122").Value = 15664.8
$ws.Range("L122").Value = 5115
$ws.Range("M122").Value = -13214.8
$ws.Range("N122").Value = -10015

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H16").Value = 1823.6666
$ws.Range("I16").Value = 1754.4706
$ws.Range("K16").Value = 1754.4706
$ws.Range("M16").Value = -1584.4706
$ws.Range("H22").Value = 1698.25
$ws.Range("I22").Value = 3499
$ws.Range("J22").Value = 1098
$ws.Range("K22").Value = 3499
$ws.Range("L22").Value = 1098
$ws.Range("M22").Value = -3204
$ws.Range("N22").Value = -1688
$ws.Range("H27").Value = 1698.25
$ws.Range("I27").Value = 3499
$ws.Range("J27").Value = 1098
$ws.Range("K27").Value = 3499
$ws.Range("L27").Value = 1098
$ws.Range("M27").Value = -3392
$ws.Range("N27").Value = -1312
$ws.Range("H61").Value = 45637.918
$ws.Range("I61").Value = 53920.65
$ws.Range("J61").Value = 4224.25
$ws.Range("K61").Value = 53920.65
$ws.Range("L61").Value = 4224.25
$ws.Range("M61").Value = -53718.65
$ws.Range("N61").Value = -4628.25
$ws.Range("H68").Value = 5500.125
$ws.Range("I68").Value = 4857.2856
$ws.Range("K68").Value = 4857.2856
$ws.Range("M68").Value = -4108.2856
$ws.Range("H71").Value = 5500.125
$ws.Range("I71").Value = 4857.2856
$ws.Range("K71").Value = 24286.428
$ws.Range("M71").Value = -20542.428
$ws.Range("H82").Value = 2744.5454
$ws.Range("I82").Value = 2355.7144
$ws.Range("J82").Value = 3425
$ws.Range("K82").Value = 2355.7144
$ws.Range("L82").Value = 3425
$ws.Range("M82").Value = -1994.7144
$ws.Range("N82").Value = -4147
$ws.Range("H85").Value = 2744.5454
$ws.Range("I85").Value = 2355.7144
$ws.Range("J85").Value = 3425
$ws.Range("K85").Value = 2355.7144
$ws.Range("L85").Value = 3425
$ws.Range("M85").Value = -1107.7144
$ws.Range("N85").Value = -5921
$ws.Range("H92").Value = 57389
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 57389
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 57389
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -62381
$ws.Range("H93").Value = 36982.332
$ws.Range("I93").Value = 4848.75
$ws.Range("J93").Value = 101249.5
$ws.Range("K93").Value = 4848.75
$ws.Range("L93").Value = 101249.5
$ws.Range("M93").Value = -3600.75
$ws.Range("N93").Value = -103745.5
$ws.Range("H100").Value = 3849629.2
$ws.Range("I100").Value = 7145896.5
$ws.Range("J100").Value = 3984.1667
$ws.Range("K100").Value = 7145896.5
$ws.Range("L100").Value = 3984.1667
$ws.Range("M100").Value = -7145355.5
$ws.Range("N100").Value = -5066.1667
$ws.Range("H113").Value = 45637.918
$ws.Range("I113").Value = 53920.65
$ws.Range("J113").Value = 4224.25
$ws.Range("K113").Value = 53920.65
$ws.Range("L113").Value = 4224.25
$ws.Range("M113").Value = -51750.65
$ws.Range("N113").Value = -8564.25
$ws.Range("H122").Value = 3872.75
$ws.Range("I122").Value = 3711
$ws.Range("K122").Value = 11133
$ws.Range("M122").Value = -8683
$ws.Range("H132").Value = 15179.637
$ws.Range("I132").Value = 16496.625
$ws.Range("J132").Value = 11667.667
$ws.Range("K132").Value = 49489.875
$ws.Range("L132").Value = 35003.001
$ws.Range("M132").Value = -46959.875
$ws.Range("N132").Value = -40063.001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H81").Value = 1631.5454
$ws.Range("I81").Value = 789.8
$ws.Range("K81").Value = 1579.6
$ws.Range("M81").Value = -518.5999999999999
$ws.Range("H84").Value = 1631.5454
$ws.Range("I84").Value = 789.8
$ws.Range("K84").Value = 7898
$ws.Range("M84").Value = -2594
$ws.Range("H122").Value = 3746.1143
$ws.Range("I122").Value = 2904.5652
$ws.Range("K122").Value = 8713.6956
$ws.Range("M122").Value = -6263.695599999999
$ws.Range("H126").Value = 5786.0557
$ws.Range("I126").Value = 3549.6924
$ws.Range("K126").Value = 10649.0772
$ws.Range("M126").Value = -8179.0772
$ws.Range("H132").Value = 4911.5405
$ws.Range("I132").Value = 3813.4285
$ws.Range("K132").Value = 11440.2855
$ws.Range("M132").Value = -8910.2855
